# =====================================================================
# feat(excel-to-json): allow comments in class and property definitions
#
# The "classes" sheet gains per-class rdfs:comment text in four
# languages (en/de/fr/it) and a couple of translations are corrected.
# =====================================================================

$wb = $excel.ActiveWorkbook
$wsClasses = $wb.Worksheets.Item("classes")

# --- widen columns C:I to fit the new language / comment columns ---
$wsClasses.Columns.Item(3).ColumnWidth = 20.5
$wsClasses.Columns.Item(4).ColumnWidth = 18.33203125
$wsClasses.Columns.Item(5).ColumnWidth = 18.33203125
$wsClasses.Columns.Item(6).ColumnWidth = 16.33203125
$wsClasses.Columns.Item(7).ColumnWidth = 15.1640625
$wsClasses.Columns.Item(8).ColumnWidth = 16.33203125
$wsClasses.Columns.Item(9).ColumnWidth = 15.1640625

# --- Header row (A1:J1): name, en, de, fr, it, comment_en, comment_de, comment_fr, comment_it, super ---
$wsClasses.Range("A1").Value = 'name'
$wsClasses.Range("B1").Value = 'en'
$wsClasses.Range("C1").Value = 'de'
$wsClasses.Range("D1").Value = 'fr'
$wsClasses.Range("E1").Value = 'it'
$wsClasses.Range("F1").Value = 'comment_en'
$wsClasses.Range("G1").Value = 'comment_de'
$wsClasses.Range("H1").Value = 'comment_fr'
$wsClasses.Range("I1").Value = 'comment_it'
$wsClasses.Range("J1").Value = 'super'

# --- Row 2: Owner ---
$wsClasses.Range("A2").Value = 'Owner'
$wsClasses.Range("B2").Value = 'Owner'
$wsClasses.Range("C2").Value = 'Eigentümer'
$wsClasses.Range("D2").Value = 'Propriétaire'
$wsClasses.Range("E2").Value = 'Proprietario'
$wsClasses.Range("F2").Value = 'A strange chance put me in possession of this journal.'
$wsClasses.Range("G2").Value = 'Ein seltsamer Zufall brachte mich in den Besitz dieses Tagebuchs.  '
$wsClasses.Range("H2").Value = 'Un étrange hasard m''a mis en possession de ce journal.  '
$wsClasses.Range("I2").Value = 'Uno strano caso mi mise in possesso di questo diario.  '
$wsClasses.Range("J2").Value = 'Resource'

# --- Row 3: Title ---
$wsClasses.Range("A3").Value = 'Title'
$wsClasses.Range("B3").Value = 'Title'
$wsClasses.Range("C3").Value = 'Titel'
$wsClasses.Range("D3").Value = 'Titre'
$wsClasses.Range("E3").Value = 'Titolo'
$wsClasses.Range("F3").Value = 'I had established myself for several months in a central city in one of our southern departments, whose shore is bathed by the Mediterranean, and I was desirous of purchasing a country place in that marvellously picturesque land. '
$wsClasses.Range("G3").Value = 'Ich hatte mich für einige Monate in einer Stadt im Zentrum eines unserer südlichen Departements niedergelassen, dessen Ufer vom Mittelmeer umspült wird, und wollte ein Grundstück in diesem wunderbar malerischen Land erwerben. '
$wsClasses.Range("H3").Value = 'Je m''étais établi depuis plusieurs mois dans une ville centrale d''un de nos départements méridionaux, dont le rivage est baigné par la Méditerranée, et je désirais acheter une maison de campagne dans cette contrée merveilleusement pittoresque. '
$wsClasses.Range("I3").Value = 'Mi ero stabilito da diversi mesi in una città centrale di uno dei nostri dipartimenti del sud, la cui riva è bagnata dal Mediterraneo, e desideravo acquistare un posto in campagna in quella terra meravigliosamente pittoresca. '
$wsClasses.Range("J3").Value = 'Resource'

# --- Row 4: GenericAnthroponym ---
$wsClasses.Range("A4").Value = 'GenericAnthroponym'
$wsClasses.Range("B4").Value = 'Generic anthroponym'
$wsClasses.Range("C4").Value = 'Allgemeines Anthroponym'
$wsClasses.Range("D4").Value = 'Anthroponyme générique'
$wsClasses.Range("E4").Value = 'Antroponimo generico'
$wsClasses.Range("F4").Value = 'I had already looked at several pieces of property when, one day, the notary, who had been giving me some necessary directions for one of my explorations, said to me:'
$wsClasses.Range("G4").Value = 'Ich hatte bereits mehrere Grundstücke besichtigt, als eines Tages der Notar, der mir die notwendigen Anweisungen für eine meiner Erkundungen gegeben hatte, zu mir sagte:'
$wsClasses.Range("H4").Value = 'J''avais déjà examiné plusieurs propriétés quand, un jour, le notaire, qui me donnait des indications nécessaires pour une de mes explorations, me dit :'
$wsClasses.Range("I4").Value = 'Avevo già visto diverse proprietà quando un giorno il notaio, che mi aveva dato alcune indicazioni necessarie per una delle mie esplorazioni, mi disse'
$wsClasses.Range("J4").Value = 'Resource'

# --- Row 5: FamilyMember ---
$wsClasses.Range("A5").Value = 'FamilyMember'
$wsClasses.Range("B5").Value = 'Family member'
$wsClasses.Range("C5").Value = 'Familienmitglied'
$wsClasses.Range("D5").Value = 'Membre de la famille'
$wsClasses.Range("E5").Value = 'Membro della famiglia'
$wsClasses.Range("F5").Value = 'I have just received notice that at about eight leagues from here, in one of the most beautiful situations in the world, neither too far nor too near to the sea, there is a country house for sale.'
$wsClasses.Range("G5").Value = 'Ich habe soeben erfahren, dass etwa acht Meilen von hier, in einer der schönsten Lagen der Welt, weder zu weit noch zu nah am Meer, ein Landhaus zum Verkauf steht.'
$wsClasses.Range("H5").Value = 'Je viens de recevoir avis qu''à huit lieues d''ici environ, dans une des plus belles situations du monde, ni trop loin ni trop près de la mer, il y a une maison de campagne à vendre.'
$wsClasses.Range("I5").Value = 'Ho appena ricevuto la notizia che a circa otto leghe da qui, in una delle situazioni più belle del mondo, né troppo lontano né troppo vicino al mare, c''è una casa di campagna in vendita.'
$wsClasses.Range("J5").Value = 'Resource'

# --- Row 6: MentionedPerson ---
$wsClasses.Range("A6").Value = 'MentionedPerson'
$wsClasses.Range("B6").Value = 'Mentioned person'
$wsClasses.Range("C6").Value = 'Erwähnte Person'
$wsClasses.Range("D6").Value = 'Personne mentionnée'
$wsClasses.Range("E6").Value = 'Persona menzionata'
$wsClasses.Range("F6").Value = 'I know nothing of it whatever; but if you would like to see it, monsieur, here are the precise directions how to find it. '
$wsClasses.Range("G6").Value = 'Ich weiß nichts davon, aber wenn Sie es sehen möchten, Monsieur, finden Sie hier die genaue Wegbeschreibung. '
$wsClasses.Range("H6").Value = 'Je n''en sais rien du tout ; mais si vous voulez la voir, monsieur, voici les indications précises pour la trouver. '
$wsClasses.Range("I6").Value = 'Non ne so nulla; ma se volete vederla, signore, eccovi le indicazioni precise per trovarla. '
$wsClasses.Range("J6").Value = 'Resource'

# --- Row 7: Alias ---
$wsClasses.Range("A7").Value = 'Alias'
$wsClasses.Range("B7").Value = 'Alias'
$wsClasses.Range("C7").Value = 'Alias'
$wsClasses.Range("D7").Value = 'Alias'
$wsClasses.Range("E7").Value = 'Alias'
$wsClasses.Range("F7").Value = 'You will have to arrange the affair with the curé of the village of ——."'
$wsClasses.Range("G7").Value = 'Sie werden die Angelegenheit mit dem Pfarrer des Dorfes -- regeln müssen."'
$wsClasses.Range("H7").Value = 'Vous devrez arranger l''affaire avec le curé du village de --."'
$wsClasses.Range("I7").Value = 'Dovrete organizzare l''affare con il curato del villaggio di --".'
$wsClasses.Range("J7").Value = 'Resource'

# --- view-state: re-select the whole sheet on Owner (matches the saved file) ---
$wsOwner = $wb.Worksheets.Item("Owner")
$wsOwner.Cells.Select()

# --- the "classes" sheet becomes the active / visible tab, zoomed to 130% ---
$wsClasses.Activate()
try { $excel.ActiveWindow.Zoom = 130 } catch { }
$wsClasses.Range("B10").Select()

